$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Add new row data (row 15): SN 11, Pending Task "DO Schedule Correction"
$ws.Range("A15").Value = 11
$ws.Range("B15").Value = "DO Schedule Correction"

# Update selection to J10
$ws.Range("J10").Select()
